$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "A branch of machine learning known as 'deep neural networks' has made..."
#   -> "A branch of machine learning known as 'deep learning' has made..."
$r1 = $d.Content
$found1 = $r1.Find.Execute("neural networks", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "learning", 2)
if (-not $found1) {
    Write-Output "WARNING: could not find 'neural networks'"
}

# --- Edit 2 -----------------------------------------------------------
# "...combining it with natural-language learning could one day lead to..."
#   -> "...combining it with natural language processing could one day lead to..."
$r2 = $d.Content
$found2 = $r2.Find.Execute("natural-language learning", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "natural language processing", 2)
if (-not $found2) {
    Write-Output "WARNING: could not find 'natural-language learning'"
}
